$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23-34 down to 24-35.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly price data.
$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C23").Value = "Los Lagos"
$ws.Range("D23").Value = 44518
$ws.Range("E23").Value = 10
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100103
$ws.Range("H23").Value = "Frutos de hueso (carozo)"
$ws.Range("I23").Value = 100103001
$ws.Range("J23").Value = "Cereza"
$ws.Range("K23").Value = "Early Burlat"
$ws.Range("L23").Value = "Segunda"
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 13000
$ws.Range("O23").Value = 13500
$ws.Range("P23").Value = 13250
$ws.Range("Q23").Value = "$/bandeja 5 kilos"
$ws.Range("R23").Value = "Provincia de Curicó"
$ws.Range("S23").Value = 2650
$ws.Range("T23").Value = 5
